$d = $word.ActiveDocument

$xml15 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="48ACB625" w14:textId="6087D153" w:rsidR="00F27CA6" w:rsidRDefault="00F27CA6"><w:pPr><w:spacing w:line="276" w:lineRule="auto"/></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$d.Paragraphs(15).Range.InsertXML($xml15)

$xml14 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="75CD7A88" w14:textId="630E1DC8" w:rsidR="003A6077" w:rsidRDefault="0091230B" w:rsidP="005D3174"><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">Author Contributions: </w:t></w:r><w:r w:rsidR="005D3174" w:rsidRPr="00737DF0"><w:t>A</w:t></w:r><w:r w:rsidR="005D3174"><w:t>CS</w:t></w:r><w:r w:rsidR="005D3174" w:rsidRPr="00737DF0"><w:t xml:space="preserve"> conceived the ideas and designed methodology; </w:t></w:r><w:r><w:t>BPME</w:t></w:r><w:r><w:t xml:space="preserve"> and ACS</w:t></w:r><w:r><w:t xml:space="preserve"> analyzed the data; A</w:t></w:r><w:r><w:t>CS</w:t></w:r><w:r><w:t xml:space="preserve"> led the writing of the manuscript. </w:t></w:r><w:r w:rsidR="005D3174"><w:t>ACS</w:t></w:r><w:r w:rsidR="005D3174" w:rsidRPr="00737DF0"><w:t xml:space="preserve"> and </w:t></w:r><w:r w:rsidR="005D3174"><w:t>BPME</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>contributed critically to the drafts and gave final approval for publication.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$d.Paragraphs(14).Range.InsertXML($xml14)

$xml13 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="6DD34008" w14:textId="7C0DC79B" w:rsidR="0091230B" w:rsidRPr="00E6736B" w:rsidRDefault="00E6736B" w:rsidP="005D3174"><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">Ethics Statement: </w:t></w:r><w:r w:rsidR="0042382D"><w:t xml:space="preserve">This research </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>was conducted</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> in compliance with t</w:t></w:r><w:r w:rsidR="0042382D" w:rsidRPr="0042382D"><w:t>he Environment and Climate Change Canada Values and Ethics Code</w:t></w:r><w:r w:rsidR="0042382D"><w:t>.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$d.Paragraphs(13).Range.InsertXML($xml13)

$xml12 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="6768A63A" w14:textId="3A91E12A" w:rsidR="00DE1D8E" w:rsidRDefault="00DE1D8E" w:rsidP="005D3174"><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">Funding Statement: </w:t></w:r><w:r w:rsidR="003C447D"><w:t xml:space="preserve">This work </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>was supported</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> by operating funds from Environment and Climate Change Canada</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$d.Paragraphs(12).Range.InsertXML($xml12)

$xml11 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="1CA4C9DC" w14:textId="77777777" w:rsidR="003C447D" w:rsidRDefault="003C447D" w:rsidP="003C447D"><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:lastRenderedPageBreak/><w:t xml:space="preserve">Data Depository: </w:t></w:r><w:r><w:t xml:space="preserve">R scripts to download the BBS data and to perform the analyses in this paper and </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>are archived</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> at </w:t></w:r><w:hyperlink r:id="rId8" w:history="1"><w:r w:rsidRPr="00410DD8"><w:rPr><w:rStyle w:val="Hyperlink"/></w:rPr><w:t>www.github.com/AdamCSmithCWS/GAM_Paper_Script</w:t></w:r></w:hyperlink></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$d.Paragraphs(11).Range.InsertXML($xml11)

$xml10 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="4F0BCAC2" w14:textId="67B5AA01" w:rsidR="005D3174" w:rsidRDefault="005D3174" w:rsidP="005D3174"><w:r><w:t xml:space="preserve">We sincerely thank the thousands of U.S. and Canadian participants who annually perform and coordinate the North American Breeding Bird Survey. We also wish to acknowledge </w:t></w:r><w:r w:rsidR="002562F4"><w:t>Courtney Amunds</w:t></w:r><w:r w:rsidR="00E8033F"><w:t>o</w:t></w:r><w:r w:rsidR="002562F4"><w:t xml:space="preserve">n for sharing some code on similar </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t xml:space="preserve">models and </w:t></w:r><w:r><w:t>John Sauer</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> and Bill Link for sharing code that helped with the cross-validations and for many spirited</w:t></w:r><w:r w:rsidR="003C447D"><w:t>,</w:t></w:r><w:r><w:t xml:space="preserve"> collegial discussions that have informed this work. We also thank the many biologists within the Canadian Wildlife Service and other users of the BBS status and trend estimates whose insightful questions and suggestions motived much of this work, including Charles Francis, Marie-Anne Hudson, Veronica Aponte, Marcel Gahbauer, Pete Blancher, and Ken Rosenberg.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$d.Paragraphs(10).Range.InsertXML($xml10)
